# Update the cryptocurrency price/volume table with the latest scraped values.
# Values in column D that look like plain numbers are prefixed with a leading
# apostrophe so Excel stores them as text (matching the original inlineStr
# cells), exactly as it would if typed in the UI.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.957.13'
$ws.Range('E2').Value = '  -0.07%  '
$ws.Range('D3').Value = '2.584.55'
$ws.Range('E3').Value = '  +1.32%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '''583.25'
$ws.Range('E5').Value = '  +0.78%  '
$ws.Range('D6').Value = '''146.83'
$ws.Range('E6').Value = '  -0.22%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  +2.09%  '
$ws.Range('E9').Value = '  +2.47%  '
$ws.Range('E10').Value = '  +2.65%  '
$ws.Range('E11').Value = '  -0.08%  '
$ws.Range('E12').Value = '  -0.30%  '
$ws.Range('D13').Value = '''27.32'
$ws.Range('E13').Value = '  +0.55%  '
$ws.Range('D14').Value = '3.047.04'
$ws.Range('E14').Value = '  +1.36%  '
$ws.Range('D15').Value = '62.802.44'
$ws.Range('E15').Value = '  -0.23%  '
$ws.Range('E16').Value = '  +3.25%  '
$ws.Range('D17').Value = '2.585.07'
$ws.Range('E17').Value = '  +1.19%  '
$ws.Range('D18').Value = '''11.32'
$ws.Range('E18').Value = '  -0.10%  '
$ws.Range('D19').Value = '''341.55'
$ws.Range('E19').Value = '  +1.80%  '
$ws.Range('D20').Value = '''4.37'
$ws.Range('E20').Value = '  +0.74%  '
$ws.Range('D21').Value = '''6.69'
$ws.Range('E21').Value = '  -0.73%  '
$ws.Range('E22').Value = '  +0.01%  '
$ws.Range('E23').Value = '  +2.42%  '
$ws.Range('D24').Value = '2.708.53'
$ws.Range('E24').Value = '  +1.12%  '
$ws.Range('E25').Value = '  -1.56%  '
$ws.Range('E26').Value = '  -0.79%  '
$ws.Range('D27').Value = '''0.999'
$ws.Range('E27').Value = '  +0.40%  '
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').Value = '''8.34'
$ws.Range('E28').Value = '  -0.26%  '
$ws.Range('B29').Value = 'Aptos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D29').Value = '''7.85'
$ws.Range('E29').Value = '  +7.29%  '
$ws.Range('D30').Value = '''1.45'
$ws.Range('E30').Value = '  -2.19%  '
$ws.Range('D31').Value = '''1.93'
$ws.Range('E31').Value = '  +1.26%  '
$ws.Range('D32').Value = '0.0₃0822'
$ws.Range('E32').Value = '  +1.17%  '
$ws.Range('D33').Value = '''467.70'
$ws.Range('E33').Value = '  +13.76%  '
$ws.Range('D34').Value = '''175.30'
$ws.Range('E34').Value = '  -1.15%  '
$ws.Range('D35').Value = '''1.61'
$ws.Range('E35').Value = '  +4.01%  '
$ws.Range('E36').Value = '  +0.07%  '
$ws.Range('E37').Value = '  +0.31%  '
$ws.Range('D38').Value = '''19.01'
$ws.Range('E38').Value = '  -0.60%  '
$ws.Range('E39').Value = '  +4.30%  '
$ws.Range('E41').Value = '  -1.90%  '
$ws.Range('D42').Value = '''157.77'
$ws.Range('E42').Value = '  +4.32%  '
$ws.Range('E43').Value = '  +0.07%  '
$ws.Range('E44').Value = '  +5.13%  '
$ws.Range('D45').Value = '''21.12'
$ws.Range('E45').Value = '  +1.55%  '
$ws.Range('D46').Value = '''0.0540'
$ws.Range('E46').Value = '  +0.55%  '
$ws.Range('E47').Value = '  -0.18%  '
$ws.Range('E48').Value = '  -0.90%  '
$ws.Range('D49').Value = '''18.42'
$ws.Range('E49').Value = '  +1.01%  '
$ws.Range('E50').Value = '  +0.57%  '
